$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Matn1"
$ws.Range("C2").Value = "Itgb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.306399
$ws.Range("H2").Value = 0.919197
$ws.Range("I2").Value = 0.3017954127735523
$ws.Range("J2").Value = 0.3017954127735522
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 164.5772705
$ws.Range("N2").Value = 329.154541
$ws.Range("O2").Value = 0.2320765473082729
$ws.Range("P2").Value = 0.1805598029509348
$ws.Range("Q2").Value = 50.42631110392951
$ws.Range("R2").Value = 302.557866623577
$ws.Range("S2").Value = 0.07003963738996105
$ws.Range("T2").Value = 0.05449212026188861

# Row 3
$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Matn1"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.306399
$ws.Range("H3").Value = 0.919197
$ws.Range("I3").Value = 0.3017954127735523
$ws.Range("J3").Value = 0.3017954127735522
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 118.764328
$ws.Range("N3").Value = 356.292984
$ws.Range("O3").Value = 0.1674740084210307
$ws.Range("P3").Value = 0.1954467673099505
$ws.Range("Q3").Value = 36.389271334872
$ws.Range("R3").Value = 327.503442013848
$ws.Range("S3").Value = 0.05054288750026632
$ws.Range("T3").Value = 0.05898493781556292

# Row 4
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Matn1"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.306399
$ws.Range("H4").Value = 0.919197
$ws.Range("I4").Value = 0.3017954127735523
$ws.Range("J4").Value = 0.3017954127735522
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 84.531957
$ws.Range("N4").Value = 253.595871
$ws.Range("O4").Value = 0.1192016653221345
$ws.Range("P4").Value = 0.1391116171686985
$ws.Range("Q4").Value = 25.900507092843
$ws.Range("R4").Value = 233.104563835587
$ws.Range("S4").Value = 0.0359745157891884
$ws.Range("T4").Value = 0.04198324792502374

# Row 5
$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Matn1"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.306399
$ws.Range("H5").Value = 0.919197
$ws.Range("I5").Value = 0.3017954127735523
$ws.Range("J5").Value = 0.3017954127735522
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 79.63570166666666
$ws.Range("N5").Value = 238.907105
$ws.Range("O5").Value = 0.1122972730628175
$ws.Range("P5").Value = 0.1310540017808179
$ws.Range("Q5").Value = 24.400299354965
$ws.Range("R5").Value = 219.602694194685
$ws.Range("S5").Value = 0.03389080187733731
$ws.Range("T5").Value = 0.03955149656306778

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Matn1"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.306399
$ws.Range("H6").Value = 0.919197
$ws.Range("I6").Value = 0.3017954127735523
$ws.Range("J6").Value = 0.3017954127735522
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 121.733284
$ws.Range("N6").Value = 365.199852
$ws.Range("O6").Value = 0.1716606440086599
$ws.Range("P6").Value = 0.2003326860218846
$ws.Range("Q6").Value = 37.298956484316
$ws.Range("R6").Value = 335.690608358844
$ws.Range("S6").Value = 0.05180639491556732
$ws.Range("T6").Value = 0.0604594856700091

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Matn1"
$ws.Range("C7").Value = "Itgb1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.306399
$ws.Range("H7").Value = 0.919197
$ws.Range("I7").Value = 0.3017954127735523
$ws.Range("J7").Value = 0.3017954127735522
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 139.9082645
$ws.Range("N7").Value = 279.816529
$ws.Range("O7").Value = 0.1972898618770847
$ws.Range("P7").Value = 0.1534951247677137
$ws.Range("Q7").Value = 42.8677523345355
$ws.Range("R7").Value = 257.206514007213
$ws.Range("S7").Value = 0.0595411753012319
$ws.Range("T7").Value = 0.04632412453800006

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Matn1"
$ws.Range("C8").Value = "Itgb1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.7088550000000001
$ws.Range("H8").Value = 2.126565
$ws.Range("I8").Value = 0.6982045872264478
$ws.Range("J8").Value = 0.6982045872264477
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 164.5772705
$ws.Range("N8").Value = 329.154541
$ws.Range("O8").Value = 0.2320765473082729
$ws.Range("P8").Value = 0.1805598029509348
$ws.Range("Q8").Value = 116.6614210802775
$ws.Range("R8").Value = 699.9685264816651
$ws.Range("S8").Value = 0.1620369099183119
$ws.Range("T8").Value = 0.1260676826890462

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Matn1"
$ws.Range("C9").Value = "Itgb1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.7088550000000001
$ws.Range("H9").Value = 2.126565
$ws.Range("I9").Value = 0.6982045872264478
$ws.Range("J9").Value = 0.6982045872264477
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 118.764328
$ws.Range("N9").Value = 356.292984
$ws.Range("O9").Value = 0.1674740084210307
$ws.Range("P9").Value = 0.1954467673099505
$ws.Range("Q9").Value = 84.18668772444
$ws.Range("R9").Value = 757.6801895199601
$ws.Range("S9").Value = 0.1169311209207644
$ws.Range("T9").Value = 0.1364618294943876

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Matn1"
$ws.Range("C10").Value = "Itgb1"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7088550000000001
$ws.Range("H10").Value = 2.126565
$ws.Range("I10").Value = 0.6982045872264478
$ws.Range("J10").Value = 0.6982045872264477
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 84.531957
$ws.Range("N10").Value = 253.595871
$ws.Range("O10").Value = 0.1192016653221345
$ws.Range("P10").Value = 0.1391116171686985
$ws.Range("Q10").Value = 59.92090037923501
$ws.Range("R10").Value = 539.2881034131151
$ws.Range("S10").Value = 0.08322714953294608
$ws.Range("T10").Value = 0.09712836924367477

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Matn1"
$ws.Range("C11").Value = "Itgb1"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.7088550000000001
$ws.Range("H11").Value = 2.126565
$ws.Range("I11").Value = 0.6982045872264478
$ws.Range("J11").Value = 0.6982045872264477
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 79.63570166666666
$ws.Range("N11").Value = 238.907105
$ws.Range("O11").Value = 0.1122972730628175
$ws.Range("P11").Value = 0.1310540017808179
$ws.Range("Q11").Value = 56.45016530492501
$ws.Range("R11").Value = 508.0514877443251
$ws.Range("S11").Value = 0.07840647118548018
$ws.Range("T11").Value = 0.0915025052177501

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Matn1"
$ws.Range("C12").Value = "Itgb1"
$ws.Range("D12").Value = "Neutro"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.7088550000000001
$ws.Range("H12").Value = 2.126565
$ws.Range("I12").Value = 0.6982045872264478
$ws.Range("J12").Value = 0.6982045872264477
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 121.733284
$ws.Range("N12").Value = 365.199852
$ws.Range("O12").Value = 0.1716606440086599
$ws.Range("P12").Value = 0.2003326860218846
$ws.Range("Q12").Value = 86.29124702982
$ws.Range("R12").Value = 776.62122326838
$ws.Range("S12").Value = 0.1198542490930926
$ws.Range("T12").Value = 0.1398732003518755

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Matn1"
$ws.Range("C13").Value = "Itgb1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.7088550000000001
$ws.Range("H13").Value = 2.126565
$ws.Range("I13").Value = 0.6982045872264478
$ws.Range("J13").Value = 0.6982045872264477
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 139.9082645
$ws.Range("N13").Value = 279.816529
$ws.Range("O13").Value = 0.1972898618770847
$ws.Range("P13").Value = 0.1534951247677137
$ws.Range("Q13").Value = 99.17467283214752
$ws.Range("R13").Value = 595.048036992885
$ws.Range("S13").Value = 0.1377486865758529
$ws.Range("T13").Value = 0.1071710002297136
